# "CreateStatSummarySheets" test case: for the "Jeremy" data sheet, add a
# fresh set of 4 statistic-summary sheets (AVERAGE / STDEVPA / MIN / MAX),
# one function per sheet, each sheet being a clone of "Jeremy" with its
# data cells replaced by a formula that pulls the corresponding value from
# "Jeremy" via the given aggregate function. This mirrors the existing
# AVERAGE (n) / STDEVPA (n) / MIN (n) / MAX (n) quadruples already present
# in the workbook, continuing the numbering sequence (... (13) -> ... (14)
# for AVERAGE/STDEVPA, ... (11) -> ... (12) for MIN/MAX).

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Jeremy")

# function name -> new sheet name, in the order they should appear
$funcs = @("AVERAGE", "STDEVPA", "MIN", "MAX")
$names = @("AVERAGE (14)", "STDEVPA (14)", "MIN (12)", "MAX (12)")

for ($i = 0; $i -lt $funcs.Length; $i++) {
    $fn = $funcs[$i]
    $nm = $names[$i]

    # Clone "Jeremy" (same layout/styles/hyperlink) and place the copy at
    # the very end of the workbook.
    $after = $wb.Worksheets.Item($wb.Worksheets.Count)
    $src.Copy($null, $after)
    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $nm

    # Replace the raw data values with formulas referencing "Jeremy",
    # wrapped in the current aggregate function - same cells Jeremy itself
    # populates (A1 header cell, the A8:D8 data row, and the B11 total).
    $new.Range("A1").Formula  = "=" + $fn + "(Jeremy!`$A`$1)"
    $new.Range("A8").Formula  = "=" + $fn + "(Jeremy!`$A`$8)"
    $new.Range("B8").Formula  = "=" + $fn + "(Jeremy!`$B`$8)"
    $new.Range("C8").Formula  = "=" + $fn + "(Jeremy!`$C`$8)"
    $new.Range("D8").Formula  = "=" + $fn + "(Jeremy!`$D`$8)"
    $new.Range("B11").Formula = "=" + $fn + "(Jeremy!`$B`$11)"
}

# The last sheet created ("MAX (12)") becomes the selected tab.
$finalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$finalSheet.Activate()
